{"js": "// Update the date heading and every worked \"two-digit \u00f7 one-digit\" answer\n// cell in the table, positionally (table values repeat, so plain text\n// search/replace is unsafe \u2014 we must address each cell by row/column).\n\n// 1) Update the date paragraph (first paragraph in the body, outside the table).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].insertText(\"2026-01-18 Sunday\", Word.InsertLocation.replace);\n\n// 2) Update the table's worked-answer cells. The table has 20 rows; only\n// rows 0, 4, 8, 12, 16 carry data (5 columns each), the rest are blank\n// spacer rows. New values below follow the exact same row/column layout.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst newValues = {\n  0: [\"99\u00f77=14, 1\", \"48\u00f73=16, 0\", \"15\u00f72=7, 1\", \"94\u00f74=23, 2\", \"47\u00f73=15, 2\"],\n  4: [\"24\u00f73=8, 0\", \"22\u00f74=5, 2\", \"47\u00f79=5, 2\", \"16\u00f74=4, 0\", \"72\u00f73=24, 0\"],\n  8: [\"79\u00f79=8, 7\", \"11\u00f74=2, 3\", \"92\u00f76=15, 2\", \"86\u00f78=10, 6\", \"35\u00f77=5, 0\"],\n  12: [\"93\u00f73=31, 0\", \"69\u00f72=34, 1\", \"75\u00f73=25, 0\", \"47\u00f73=15, 2\", \"33\u00f74=8, 1\"],\n  16: [\"39\u00f73=13, 0\", \"34\u00f76=5, 4\", \"88\u00f77=12, 4\", \"51\u00f73=17, 0\", \"95\u00f76=15, 5\"],\n};\n\nfor (const rowIndex of Object.keys(newValues)) {\n  const row = parseInt(rowIndex, 10);\n  const values = newValues[rowIndex];\n  for (let col = 0; col < values.length; col++) {\n    const cell = table.getCell(row, col);\n    cell.value = values[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and every worked \"two-digit \u00f7 one-digit\" answer\n# cell in the table, positionally (several cell values repeat verbatim, so\n# a blind Find/Replace-All would clobber the wrong instance \u2014 we must\n# address each paragraph/cell by its position instead).\n\n$d = $word.ActiveDocument\n\n# 1) Update the date paragraph (first paragraph in the body, outside the table).\n$d.Paragraphs(1).Range.Text = \"2026-01-18 Sunday\"\n\n# 2) Update the table's worked-answer cells. The table has 20 rows; only\n# rows 1, 5, 9, 13, 17 (1-indexed) carry data (5 columns each) \u2014 the rest\n# are blank spacer rows. New values below follow the exact same\n# row/column layout as the source table.\n$t = $d.Tables(1)\n\n$newValues = @{\n    1  = @(\"99\u00f77=14, 1\", \"48\u00f73=16, 0\", \"15\u00f72=7, 1\", \"94\u00f74=23, 2\", \"47\u00f73=15, 2\")\n    5  = @(\"24\u00f73=8, 0\", \"22\u00f74=5, 2\", \"47\u00f79=5, 2\", \"16\u00f74=4, 0\", \"72\u00f73=24, 0\")\n    9  = @(\"79\u00f79=8, 7\", \"11\u00f74=2, 3\", \"92\u00f76=15, 2\", \"86\u00f78=10, 6\", \"35\u00f77=5, 0\")\n    13 = @(\"93\u00f73=31, 0\", \"69\u00f72=34, 1\", \"75\u00f73=25, 0\", \"47\u00f73=15, 2\", \"33\u00f74=8, 1\")\n    17 = @(\"39\u00f73=13, 0\", \"34\u00f76=5, 4\", \"88\u00f77=12, 4\", \"51\u00f73=17, 0\", \"95\u00f76=15, 5\")\n}\n\nforeach ($rowIndex in $newValues.Keys) {\n    $values = $newValues[$rowIndex]\n    for ($col = 1; $col -le $values.Length; $col++) {\n        $t.Cell($rowIndex, $col).Range.Text = $values[$col - 1]\n    }\n}\n"}
